$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(296).Insert()

$ws.Cells.Item(296, 1).Value = 9
$ws.Cells.Item(296, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(296, 3).Value = "Metropolitana"
$ws.Cells.Item(296, 4).Value = 44642
$ws.Cells.Item(296, 5).Value = 13
$ws.Cells.Item(296, 6).Value = 100112039
$ws.Cells.Item(296, 7).Value = "Ciboulette"
$ws.Cells.Item(296, 8).Value = "Sin especificar"
$ws.Cells.Item(296, 9).Value = "Primera"
$ws.Cells.Item(296, 10).Value = 160
$ws.Cells.Item(296, 11).Value = 1200
$ws.Cells.Item(296, 12).Value = 1500
$ws.Cells.Item(296, 13).Value = 1350
$ws.Cells.Item(296, 14).Value = '$/docena de atados'
$ws.Cells.Item(296, 15).Value = "Región Metropolitana"
$ws.Cells.Item(296, 16).Value = 450
$ws.Cells.Item(296, 17).Value = 3
$ws.Cells.Item(296, 18).Value = "Hortaliza"
